# Applies the "Updated cryptos list" data refresh (GitHub Actions run of
# Wed Apr 24 04:54:44 UTC 2024) to the crypto price/volume table on Sheet1.
#
# Price (column D) and 1h-volume-change (column E) values are refreshed for
# every coin row; rows 32/33 additionally swap places (Stacks now ranks above
# EthereumClassic), so their Coin/Link/Price/Volume cells are rewritten too.
#
# All of these columns hold plain text in the workbook (not numbers), so for
# any new Price value that reads like a bare number (e.g. '0.997') we prefix
# it with a leading apostrophe -- Excel's normal force-text quote-prefix --
# so the cell keeps storing/displaying it as text instead of silently
# converting it to a numeric cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '66.982.64'
$ws.Range('E2').Value = '  +1.06%  '

# Row 3
$ws.Range('D3').Value = '3.266.47'
$ws.Range('E3').Value = '  +2.80%  '

# Row 4
$ws.Range('D4').Value = '''0.997'
$ws.Range('E4').Value = '  -0.35%  '

# Row 5
$ws.Range('D5').Value = '''606.87'
$ws.Range('E5').Value = '  +1.11%  '

# Row 6
$ws.Range('D6').Value = '''158.94'
$ws.Range('E6').Value = '  +2.71%  '

# Row 7
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.15%  '

# Row 8
$ws.Range('D8').Value = '3.268.43'
$ws.Range('E8').Value = '  +2.96%  '

# Row 9
$ws.Range('D9').Value = '''0.550'
$ws.Range('E9').Value = '  +0.73%  '

# Row 10
$ws.Range('D10').Value = '''0.163'
$ws.Range('E10').Value = '  +3.39%  '

# Row 11
$ws.Range('D11').Value = '''5.87'
$ws.Range('E11').Value = '  +2.12%  '

# Row 12
$ws.Range('D12').Value = '''0.507'
$ws.Range('E12').Value = '  +0.02%  '

# Row 13
$ws.Range('D13').Value = '''0.0000273'
$ws.Range('E13').Value = '  +4.17%  '

# Row 14
$ws.Range('D14').Value = '''39.68'
$ws.Range('E14').Value = '  +2.69%  '

# Row 15
$ws.Range('D15').Value = '3.776.98'
$ws.Range('E15').Value = '  +2.00%  '

# Row 16
$ws.Range('D16').Value = '66.852.46'
$ws.Range('E16').Value = '  +0.71%  '

# Row 17
$ws.Range('D17').Value = '''7.39'
$ws.Range('E17').Value = '  +0.36%  '

# Row 18
$ws.Range('D18').Value = '3.256.08'
$ws.Range('E18').Value = '  +2.01%  '

# Row 19
$ws.Range('E19').Value = '  +1.89%  '

# Row 20
$ws.Range('D20').Value = '''511.95'
$ws.Range('E20').Value = '  +0.42%  '

# Row 21
$ws.Range('D21').Value = '''15.49'
$ws.Range('E21').Value = '  +0.68%  '

# Row 22
$ws.Range('D22').Value = '''0.757'
$ws.Range('E22').Value = '  +3.90%  '

# Row 23
$ws.Range('D23').Value = '''8.15'
$ws.Range('E23').Value = '  +0.87%  '

# Row 24
$ws.Range('D24').Value = '''14.83'
$ws.Range('E24').Value = '  +0.08%  '

# Row 25
$ws.Range('D25').Value = '''86.36'
$ws.Range('E25').Value = '  +2.45%  '

# Row 26
$ws.Range('D26').Value = '''0.175'
$ws.Range('E26').Value = '  +96.54%  '

# Row 27
$ws.Range('E27').Value = '  +0.33%  '

# Row 28
$ws.Range('D28').Value = '''3.03'
$ws.Range('E28').Value = '  +1.15%  '

# Row 29
$ws.Range('D29').Value = '''9.16'
$ws.Range('E29').Value = '  +0.17%  '

# Row 30
$ws.Range('D30').Value = '''2.38'
$ws.Range('E30').Value = '  -0.06%  '

# Row 31
$ws.Range('D31').Value = '''6.98'
$ws.Range('E31').Value = '  -0.81%  '

# Row 32
$ws.Range('B32').Value = 'Stacks'
$ws.Range('C32').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D32').Value = '''2.89'
$ws.Range('E32').Value = '  -5.42%  '

# Row 33
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '''28.45'
$ws.Range('E33').Value = '  +1.80%  '

# Row 34
$ws.Range('D34').Value = '''0.998'
$ws.Range('E34').Value = '  -0.32%  '

# Row 35
$ws.Range('D35').Value = '''1.17'
$ws.Range('E35').Value = '  -2.97%  '

# Row 36
$ws.Range('D36').Value = '''6.44'
$ws.Range('E36').Value = '  -1.03%  '

# Row 37
$ws.Range('D37').Value = '0.0₃0804'
$ws.Range('E37').Value = '  +20.29%  '

# Row 38
$ws.Range('D38').Value = '''3.37'
$ws.Range('E38').Value = '  +19.79%  '

# Row 39
$ws.Range('D39').Value = '''55.51'
$ws.Range('E39').Value = '  +1.51%  '

# Row 40
$ws.Range('D40').Value = '''496.34'
$ws.Range('E40').Value = '  -2.56%  '

# Row 41
$ws.Range('D41').Value = '''0.0429'
$ws.Range('E41').Value = '  +1.99%  '

# Row 42
$ws.Range('E42').Value = '  +2.91%  '

# Row 43
$ws.Range('D43').Value = '''8.82'
$ws.Range('E43').Value = '  +0.14%  '

# Row 44
$ws.Range('D44').Value = '''0.297'
$ws.Range('E44').Value = '  -0.63%  '

# Row 45
$ws.Range('D45').Value = '''2.49'
$ws.Range('E45').Value = '  +2.84%  '

# Row 46
$ws.Range('D46').Value = '2.957.93'
$ws.Range('E46').Value = '  +4.11%  '

# Row 47
$ws.Range('D47').Value = '''28.65'
$ws.Range('E47').Value = '  +1.50%  '

# Row 48
$ws.Range('D48').Value = '''2.43'
$ws.Range('E48').Value = '  +0.97%  '

# Row 49
$ws.Range('D49').Value = '''0.120'
$ws.Range('E49').Value = '  +2.90%  '

# Row 50
$ws.Range('E50').Value = '  -0.02%  '

# Row 51
$ws.Range('D51').Value = '''2.57'
$ws.Range('E51').Value = '  +1.81%  '
